$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.046.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.820.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "701.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +8.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.19"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.818.66"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.25"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.89%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.459.72"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.812.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.952.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.85"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.21"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +17.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "480.92"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.68%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.17"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.48"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.966.98"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.11"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +15.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.31"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.55"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.185"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.59"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.26"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.04"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +12.94%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.982"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000324"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +22.18%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.20%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.85"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.68"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.301"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.57"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.77%  "